$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Apr 2020 to May 2020"
$ws.Range("A3").Select()
